$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The script finished running: the "SCRIPT RESULT" column (I) now matches
# the "RESULT EXCEL" column (H) for rows 18-20, so their STATUS (column J)
# flips from NOK to OK.
$ws.Range("I18").Value = $ws.Range("H18").Value2
$ws.Range("J18").Value = "OK"

$ws.Range("I19").Value = $ws.Range("H19").Value2
$ws.Range("J19").Value = "OK"

$ws.Range("I20").Value = $ws.Range("H20").Value2
$ws.Range("J20").Value = "OK"

# Two stray formatted-but-empty cells appear (artifacts of selecting /
# formatting a range in the UI). Build the required formatting on a
# scratch cell, then copy it onto the target cells, then clean up.
$scratch = $ws.Range("ZZ9999")

# O16 picks up the same "blank helper" look already used at O7 / N9 / N24.
$ws.Range("N9").Copy() | Out-Null
$ws.Range("O16").PasteSpecial(-4122) | Out-Null

# M23 picks up a new, distinct blank style.
$scratch.Font.Name = "Calibri"
$scratch.Copy() | Out-Null
$ws.Range("M23").PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

# The workbook was left scrolled to/selecting M23.
$ws.Range("M23").Select() | Out-Null

Write-Host "done"
